$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (C4)
$ws.Range("B6").Value = 78.799999999999997
$ws.Range("C6").Value = -38.950000000000003

# Row 7 (C5)
$ws.Range("B7").Value = 77.290000000000006
$ws.Range("C7").Value = -42.280000000000001

# Row 8 (C6)
$ws.Range("B8").Value = 80.590000000000003
$ws.Range("C8").Value = -42.280000000000001

# Row 24 (L1)
$ws.Range("B24").Value = 78.939999999999998
$ws.Range("C24").Value = -42.280000000000001

# Row 34 (RN1)
$ws.Range("B34").Value = 63.68
$ws.Range("C34").Value = -49.770000000000003
